# Adds the "04-Nov-2025" progress update:
#   - Training Dashboard: PERIOD TO EXPIRE (col H) decreases by 1 day and
#     LAST UPDATE (col I) moves from 03-Nov-2025 to 04-Nov-2025 for rows 3-17.
#   - Exam Dashboard: a new "Cs Hoist" low-percentage result row is inserted
#     above the TOTAL AVERAGE row, and the total average is recalculated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Training Dashboard
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Dashboard")

$periodUpdates = @{
    3  = 623
    4  = 625
    5  = 623
    6  = 625
    7  = 623
    8  = 624
    9  = 625
    10 = 624
    11 = 625
    12 = 626
    13 = 626
    14 = 626
    15 = 259
    16 = 352
    17 = 299
}

# Force the LAST UPDATE column to be treated as text so Excel does not
# silently turn the "dd-Mmm-yyyy" strings into date serial numbers.
$ws1.Range("I3:I17").NumberFormat = "@"

foreach ($row in $periodUpdates.Keys) {
    $ws1.Range("H$row").Value = $periodUpdates[$row]
    $ws1.Range("I$row").Value = "04-Nov-2025"

    # Re-apply the row's normal style (General number format, same border/
    # fill/font) by pulling it from column J of the same row, which is left
    # untouched by this edit, so the temporary text number format above
    # does not linger as a distinct cell style.
    $ws1.Range("J$row").Copy()
    $ws1.Range("I$row").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet 2: Exam Dashboard
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# Insert a new row above the TOTAL AVERAGE row (currently row 6); this
# pushes TOTAL AVERAGE down to row 7.
$ws2.Range("A6:G6").Insert()

# Copy the formatting of the row above (an existing "low percentage" exam
# result) onto the freshly inserted row so it gets the same pink style.
$ws2.Range("A5:G5").Copy()
$ws2.Range("A6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The exam date and marks columns must stay as plain text (like all the
# other rows), not get auto-converted to a date serial / percentage number.
$ws2.Range("C6:D6").NumberFormat = "@"

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "Cs Hoist"
$ws2.Range("C6").Value = "30-Oct-2025"
$ws2.Range("D6").Value = "40.74%"
$ws2.Range("E6").Value = "low percentage"
$ws2.Range("F6").Value = "This is a low mark, please retake the exam and improve your score. date is valid"

# Re-apply the pink row style on top of the values so the temporary text
# format above collapses back onto the shared "low percentage" style.
$ws2.Range("A5:G5").Copy()
$ws2.Range("A6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the TOTAL AVERAGE row, which has shifted down to row 7.
$ws2.Range("D7").NumberFormat = "@"
$ws2.Range("D7").Value = "48.94%"

# Restore D7's normal style from C7 (same row, untouched) so it keeps the
# plain General-format TOTAL AVERAGE row style instead of a forked one.
$ws2.Range("C7").Copy()
$ws2.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
